# Apply cell updates from the cryptos list refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRef, $value)
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "58.988.68"
Set-TextValue "E2" "  -2.36%  "
Set-TextValue "D3" "2.661.11"
Set-TextValue "E3" "  -0.90%  "
Set-TextValue "E4" "  -0.04%  "
Set-TextValue "D5" "525.41"
Set-TextValue "E5" "  +0.52%  "
Set-TextValue "D6" "144.40"
Set-TextValue "E6" "  -1.20%  "
Set-TextValue "E7" "  +0.24%  "
Set-TextValue "D8" "0.569"
Set-TextValue "E8" "  -1.02%  "
Set-TextValue "D9" "7.00"
Set-TextValue "E9" "  +8.27%  "
Set-TextValue "E10" "  -2.31%  "
Set-TextValue "E11" "  -2.04%  "
Set-TextValue "D12" "0.131"
Set-TextValue "E12" "  +1.44%  "
Set-TextValue "D13" "3.130.22"
Set-TextValue "E13" "  -0.92%  "
Set-TextValue "D14" "58.990.03"
Set-TextValue "E14" "  -2.39%  "
Set-TextValue "E15" "  -0.98%  "
Set-TextValue "B16" "ShibaInu"
Set-TextValue "C16" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D16" "0.0000137"
Set-TextValue "E16" "  -1.57%  "
Set-TextValue "B17" "WrappedEther"
Set-TextValue "C17" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D17" "2.646.71"
Set-TextValue "E17" "  -4.29%  "
Set-TextValue "D18" "338.87"
Set-TextValue "E18" "  -3.39%  "
Set-TextValue "E19" "  -3.21%  "
Set-TextValue "D20" "10.38"
Set-TextValue "E20" "  -2.16%  "
Set-TextValue "D21" "6.40"
Set-TextValue "E21" "  +1.33%  "
Set-TextValue "D22" "1.01"
Set-TextValue "E22" "  +0.74%  "
Set-TextValue "D23" "64.28"
Set-TextValue "E23" "  +2.22%  "
Set-TextValue "E24" "  -0.93%  "
Set-TextValue "E25" "  -1.83%  "
Set-TextValue "D26" "0.999"
Set-TextValue "E26" "  +0.48%  "
Set-TextValue "D27" "0.0₃0801"
Set-TextValue "E27" "  -1.68%  "
Set-TextValue "D28" "7.10"
Set-TextValue "E28" "  -2.36%  "
Set-TextValue "D29" "6.69"
Set-TextValue "E29" "  -2.65%  "
Set-TextValue "E30" "  +0.12%  "
Set-TextValue "E31" "  -0.16%  "
Set-TextValue "D32" "18.86"
Set-TextValue "E32" "  -1.09%  "
Set-TextValue "D33" "150.65"
Set-TextValue "E33" "  +1.48%  "
Set-TextValue "E34" "  -4.10%  "
Set-TextValue "E35" "  -3.93%  "
Set-TextValue "D36" "0.894"
Set-TextValue "E36" "  -5.81%  "
Set-TextValue "D37" "0.875"
Set-TextValue "E37" "  -0.40%  "
Set-TextValue "D38" "36.87"
Set-TextValue "E38" "  +0.03%  "
Set-TextValue "E39" "  -5.72%  "
Set-TextValue "D40" "3.59"
Set-TextValue "E40" "  -2.96%  "
Set-TextValue "E41" "  +0.45%  "
Set-TextValue "E42" "  +0.27%  "
Set-TextValue "D43" "275.96"
Set-TextValue "E43" "  -2.24%  "
Set-TextValue "D44" "19.89"
Set-TextValue "E44" "  -0.31%  "
Set-TextValue "E45" "  -1.71%  "
Set-TextValue "D46" "10.66"
Set-TextValue "E46" "  +1.98%  "
Set-TextValue "D47" "2.051.63"
Set-TextValue "E47" "  -3.57%  "
Set-TextValue "E48" "  -1.67%  "
Set-TextValue "E49" "  -3.36%  "
Set-TextValue "E50" "  -2.29%  "
Set-TextValue "D51" "18.96"
Set-TextValue "E51" "  -1.24%  "
